$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 'Favourites'
$ws.Range("H5").Value = 'Implement ''RemoveFavourite'' from Details & Index'
$ws.Range("G6").Value = $null
$ws.Range("H6").Value = $null
$ws.Range("G7").Value = 'Login'
$ws.Range("H7").Value = 'User can log in from modal'
$ws.Range("G8").Value = $null
$ws.Range("H8").Value = $null
$ws.Range("H9").Value = 'Finish update - complete text, fa icons etc'
$ws.Range("G10").Value = 'Homepage'
$ws.Range("H10").Value = 'Add Google chart to Homepage'
$ws.Range("G11").Value = $null
$ws.Range("H11").Value = $null
$ws.Range("G12").Value = 'Reviews'
$ws.Range("H12").Value = 'Reviews - and Ratings - Stars fa icons'
$ws.Range("G13").Value = $null
$ws.Range("H13").Value = $null
$ws.Range("H14").Value = 'Add Bool for Featured Cars, Cars On special '
$ws.Range("G15").Value = 'Special Offer'
$ws.Range("H15").Value = 'link to special offer/ featuredlist on click Call To Action Links'
$ws.Range("H16").Value = $null
$ws.Range("H17").Value = 'Finish Categories Sidebar'
$ws.Range("H18").Value = 'Owl Carousel data-img-zoom not working'
$ws.Range("H19").Value = 'Check Footer links'
$ws.Range("H20").Value = $null
$ws.Range("H21").Value = 'Return IHttpActionResult from apiControllers  - Ok, Bad Result etc.'
$ws.Range("H22").Value = 'Unit Tests - Nunit'
$ws.Range("H23").Value = 'Jasmine Tests JS '
$ws.Range("H24").Value = 'AutoComplete on text boxes???'
$ws.Range("G25").Value = $null
$ws.Range("H25").Value = $null
$ws.Range("G26").Value = 'Validation'
$ws.Range("H26").Value = 'Validation logic and add modelstate errors etc.'
$ws.Range("G27").Value = $null
$ws.Range("H27").Value = $null
$ws.Range("G28").Value = 'Pagination'
$ws.Range("H28").Value = 'Go direct to certain page - eg pg 17'
$ws.Range("G29").Value = $null
$ws.Range("H29").Value = $null
$ws.Range("G30").Value = 'Error Logging'
$ws.Range("H30").Value = 'Log errors with NLog?'
$ws.Range("G31").Value = $null
$ws.Range("H31").Value = $null
$ws.Range("J31").Value = 'Bug: Add recommended cars to cookie'
$ws.Range("G32").Value = 'Admin Area'
$ws.Range("H32").Value = 'Update cars etc.'
$ws.Range("J32").Value = 'Refactor: Run everything off of "Search"'
$ws.Range("G33").Value = $null
$ws.Range("H33").Value = $null
$ws.Range("J33").Value = $null
$ws.Range("H34").Value = 'Send email on registration'
$ws.Range("J34").Value = 'Cost per day per car - different amounts e.g. Category A, B , C'
$ws.Range("G35").Value = 'Email'
$ws.Range("H35").Value = 'Send email on mailing list sign up'
$ws.Range("G36").Value = $null
$ws.Range("H36").Value = $null
$ws.Range("J36").Value = 'Mega Menu Links'
$ws.Range("G37").Value = 'Recommended'
$ws.Range("H37").Value = 'Recommended cars more tailored to users favourite cars'
$ws.Range("G38").Value = $null
$ws.Range("H38").Value = $null
$ws.Range("G39").Value = 'Grid View'
$ws.Range("H39").Value = 'Add Car List as Grid View'
$ws.Range("G41").Value = 'Cookies'
$ws.Range("H41").Value = 'Change from email address to ASPnet_Users table ID'

$ws.Range("H25").Select()
